# Modelagem testes de calculo de dano
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# --- Remove the extra, unused sheet ---
$wb.Worksheets.Item("Plan3").Delete()

# --- Column widths: A wider to fit new labels, B gets an explicit width ---
$ws.Columns.Item(1).ColumnWidth = 23.85546875
$ws.Columns.Item(2).ColumnWidth = 9.140625

# --- New "Weapon Calculator" section ---
$ws.Range("A19").Value = "Weapon Calculator"

$ws.Range("A20").Value = "Lex Prime"
$ws.Range("B20").Formula = "=SUM(B21:B23)"

$ws.Range("A21").Value = "Impact"
$ws.Range("B21").Value = 15

$ws.Range("A22").Value = "Puncture"
$ws.Range("B22").Value = 120

$ws.Range("A23").Value = "Slash"
$ws.Range("B23").Value = 15

$ws.Range("A24").Value = "Heat (Primed HC +165%)"
$ws.Range("B24").Formula = "=B20*2.65"

$ws.Range("A25").Value = "Total Damage"
$ws.Range("B25").Formula = "=SUM(B21:B24)"

$ws.Range("A27").Value = "Xata Whisper (+17%)"
$ws.Range("B27").Formula = "=B25*1.17"

$ws.Range("A28").Value = "Impact"
$ws.Range("B28").Value = 15

$ws.Range("A29").Value = "Puncture"
$ws.Range("B29").Value = 120

$ws.Range("A30").Value = "Slash"
$ws.Range("B30").Value = 15

$ws.Range("A31").Value = "Heat"
$ws.Range("B31").Formula = "=B24"

$ws.Range("A32").Value = "Void"
$ws.Range("B32").Formula = "=B26*0.17"

$ws.Range("A34").Value = "Roar (+50%)"
$ws.Range("B34").Formula = "=SUM(B35:B38)"

$ws.Range("A35").Value = "Impact"
$ws.Range("B35").Formula = "=B21*1.5"

$ws.Range("A36").Value = "Puncture"
$ws.Range("B36").Formula = "=B22*1.5"

$ws.Range("A37").Value = "Slash"
$ws.Range("B37").Formula = "=B23*1.5"

$ws.Range("A38").Value = "Heat"
$ws.Range("B38").Formula = "=B24*1.5"

$ws.Range("A40").Value = "Vex Armor (+275%)"
$ws.Range("B40").Formula = "=SUM(B41:B44)"

$ws.Range("A41").Value = "Impact"
$ws.Range("B41").Formula = "=B21*3.75"

$ws.Range("A42").Value = "Puncture"
$ws.Range("B42").Formula = "=B22*3.75"

$ws.Range("A43").Value = "Slash"
$ws.Range("B43").Formula = "=B23*3.75"

$ws.Range("A44").Value = "Heat"
$ws.Range("B44").Formula = "=B24"

# --- Selection / scroll position ---
$ws.Range("B19").Select()
